# Updates crypto price/volume figures per upstream symbol-list refresh.
# Each target cell holds its value as TEXT in the source workbook (t="inlineStr"),
# e.g. D2 = "310.25" and E2 = "-3.37%" — not numbers/percentages. Assigning a
# bare numeric- or percent-looking string via COM Range.Value lets Excel's
# input parser reinterpret it as a Number/Percentage cell, which would change
# the cell's type. To keep these as plain text (matching the source), we
# write the value with a leading apostrophe (forces text entry, like typing
# it in the Excel UI) and then reset the cell style to "Normal" so no stray
# quote-prefix formatting/style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "309.81"
Set-TextValue "E2" "-3.70%"
Set-TextValue "D3" "53.57"
Set-TextValue "E3" "8.04%"
Set-TextValue "D4" "5.137"
Set-TextValue "E4" "-2.65%"
Set-TextValue "D5" "0.07841"
Set-TextValue "E5" "-2.52%"
Set-TextValue "D6" "4.522"
Set-TextValue "E6" "-1.25%"
Set-TextValue "E7" "-3.90%"
Set-TextValue "D8" "1.584"
Set-TextValue "E8" "-3.74%"
Set-TextValue "D9" "0.1219"
Set-TextValue "E9" "-7.15%"
Set-TextValue "D10" "0.2027"
Set-TextValue "E10" "2.88%"
Set-TextValue "D11" "0.04721"
Set-TextValue "E11" "0.57%"
Set-TextValue "D12" "0.09459"
Set-TextValue "E12" "0.55%"
Set-TextValue "D13" "0.1043"
Set-TextValue "E13" "-0.31%"
Set-TextValue "D14" "0.001264"
Set-TextValue "E14" "-5.81%"
Set-TextValue "E15" "-0.89%"
Set-TextValue "E16" "2,031.18%"
Set-TextValue "D17" "3.337"
Set-TextValue "E17" "-0.37%"
Set-TextValue "E18" "-0.95%"
Set-TextValue "D19" "0.3418"
Set-TextValue "E19" "-1.39%"
Set-TextValue "D20" "7.982"
Set-TextValue "E20" "-1.78%"
Set-TextValue "E21" "-1.21%"
Set-TextValue "D23" "0.04131"
Set-TextValue "E23" "-0.44%"
Set-TextValue "E24" "-3.62%"
Set-TextValue "D25" "0.003954"
Set-TextValue "E25" "-8.53%"
Set-TextValue "D26" "0.0001347"
Set-TextValue "E26" "0.50%"
Set-TextValue "E38" "-3.27%"
Set-TextValue "D39" "0.05950"
Set-TextValue "E39" "-1.18%"
Set-TextValue "E40" "0.95%"
Set-TextValue "D41" "0.007904"
Set-TextValue "E41" "-1.07%"
Set-TextValue "D42" "0.1428"
Set-TextValue "E42" "-1.99%"
Set-TextValue "D43" "0.008200"
Set-TextValue "E43" "6.19%"
Set-TextValue "D44" "0.008468"
Set-TextValue "E44" "-1.83%"
Set-TextValue "D45" "0.3126"
Set-TextValue "E45" "-2.52%"
Set-TextValue "D46" "0.00007231"
Set-TextValue "E46" "9.44%"
Set-TextValue "D47" "0.00000000748"
Set-TextValue "E47" "0.49%"
Set-TextValue "D48" "0.05595"
Set-TextValue "E48" "1.70%"
Set-TextValue "D49" "0.002613"
Set-TextValue "E49" "-34.18%"
Set-TextValue "D50" "0.00002095"
Set-TextValue "E50" "0.49%"
Set-TextValue "D51" "0.0001995"
Set-TextValue "E51" "0.49%"
